# Applies the cryptocurrency price/volume refresh described in the commit diff.
# - Updates Price (column D) and Volume(1h) (column E) text values for many rows.
# - Fixes two rows whose coin data had been swapped (rows 29/30 and 46/47),
#   restoring the correct Coin / Link / Price / Volume(1h) values for each.
#
# Note: several "Price" values look like plain numbers (e.g. "101.30", "10.70").
# Excel would normally auto-convert such text to a numeric value (dropping the
# trailing zero, etc.), so for those cells we use Excel's classic leading
# apostrophe trick to force the value to be stored as text, exactly as the
# source data has it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.115.04"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "2.944.74"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'374.96"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").Value = "'101.30"
$ws.Range("E6").Value = "  -3.02%  "
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").Value = "'36.33"
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("D13").Value = "3.404.94"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "'18.07"
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "2.955.24"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "'0.993"
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("D18").Value = "'10.70"
$ws.Range("E18").Value = "  +44.07%  "
$ws.Range("D19").Value = "50.990.03"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("E20").Value = "  -6.52%  "
$ws.Range("D21").Value = "'12.44"
$ws.Range("E21").Value = "  -4.14%  "
$ws.Range("D22").Value = "0.0₃0958"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").Value = "'266.32"
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("D24").Value = "'68.68"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("E25").Value = "  +8.50%  "
$ws.Range("E26").Value = "  -2.84%  "
$ws.Range("D27").Value = "'7.61"
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.164"
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'25.64"
$ws.Range("E30").Value = "  -1.35%  "
$ws.Range("E31").Value = "  -5.71%  "
$ws.Range("D32").Value = "'10.02"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").Value = "'50.68"
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("E35").Value = "  -4.90%  "
$ws.Range("E36").Value = "  -2.32%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  +4.11%  "
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("D40").Value = "'16.27"
$ws.Range("E40").Value = "  -5.01%  "
$ws.Range("E41").Value = "  -2.84%  "
$ws.Range("D42").Value = "'2.48"
$ws.Range("E42").Value = "  -3.89%  "
$ws.Range("D43").Value = "'120.40"
$ws.Range("E43").Value = "  -4.04%  "
$ws.Range("D44").Value = "'21.44"
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.273"
$ws.Range("E46").Value = "  -3.54%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'3.33"
$ws.Range("E47").Value = "  +2.23%  "
$ws.Range("D48").Value = "'2.31"
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("D49").Value = "1.997.44"
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("D50").Value = "'0.0324"
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("E51").Value = "  +1.09%  "
